# Apply the "cryptos list" update (GitHub Actions scheduled refresh).
# Only the Price (D) and Volume(1h) (E) columns change for most rows;
# rows 47/48 (Stellar/ThetaToken) also swap position/rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.739.35'
$ws.Range('E2').Value = '  -1.53%  '
# Row 3
$ws.Range('D3').Value = '3.332.39'
$ws.Range('E3').Value = '  -1.26%  '
# Row 4
$ws.Range('E4').Value = '  +0.02%  '
# Row 5
$ws.Range('D5').Value = '''582.49'
$ws.Range('E5').Value = '  -2.07%  '
# Row 6
$ws.Range('D6').Value = '''175.74'
$ws.Range('E6').Value = '  -5.72%  '
# Row 7
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  +0.15%  '
# Row 8
$ws.Range('D8').Value = '''0.589'
$ws.Range('E8').Value = '  -1.65%  '
# Row 9
$ws.Range('D9').Value = '3.328.99'
$ws.Range('E9').Value = '  -0.78%  '
# Row 10
$ws.Range('E10').Value = '  -2.87%  '
# Row 11
$ws.Range('D11').Value = '''0.576'
$ws.Range('E11').Value = '  -1.80%  '
# Row 12
$ws.Range('D12').Value = '''45.56'
$ws.Range('E12').Value = '  -3.65%  '
# Row 13
$ws.Range('D13').Value = '''0.0000269'
$ws.Range('E13').Value = '  -3.95%  '
# Row 14
$ws.Range('D14').Value = '''658.41'
$ws.Range('E14').Value = '  +2.75%  '
# Row 15
$ws.Range('D15').Value = '3.870.84'
$ws.Range('E15').Value = '  -1.13%  '
# Row 16
$ws.Range('D16').Value = '''8.40'
$ws.Range('E16').Value = '  -1.65%  '
# Row 17
$ws.Range('D17').Value = '67.943.16'
$ws.Range('E17').Value = '  -1.50%  '
# Row 19
$ws.Range('D19').Value = '3.334.98'
$ws.Range('E19').Value = '  -1.32%  '
# Row 20
$ws.Range('D20').Value = '''17.41'
# Row 21
$ws.Range('D21').Value = '''10.94'
$ws.Range('E21').Value = '  -1.40%  '
# Row 22
$ws.Range('D22').Value = '''0.889'
$ws.Range('E22').Value = '  -2.44%  '
# Row 23
$ws.Range('D23').Value = '''5.42'
$ws.Range('E23').Value = '  +6.63%  '
# Row 24
$ws.Range('D24').Value = '''17.04'
$ws.Range('E24').Value = '  -5.09%  '
# Row 25
$ws.Range('D25').Value = '''99.53'
$ws.Range('E25').Value = '  +0.39%  '
# Row 26
$ws.Range('E26').Value = '  -5.99%  '
# Row 27
$ws.Range('D27').Value = '''2.67'
$ws.Range('E27').Value = '  -6.18%  '
# Row 28
$ws.Range('D28').Value = '''9.26'
$ws.Range('E28').Value = '  -5.59%  '
# Row 29
$ws.Range('D29').Value = '''33.50'
$ws.Range('E29').Value = '  +1.65%  '
# Row 30
$ws.Range('D30').Value = '''7.41'
$ws.Range('E30').Value = '  +8.72%  '
# Row 31
$ws.Range('D31').Value = '''8.43'
$ws.Range('E31').Value = '  -3.02%  '
# Row 32
$ws.Range('D32').Value = '''590.95'
$ws.Range('E32').Value = '  -3.24%  '
# Row 33
$ws.Range('E33').Value = '  -1.20%  '
# Row 34
$ws.Range('D34').Value = '''0.104'
$ws.Range('E34').Value = '  -1.12%  '
# Row 35
$ws.Range('D35').Value = '3.718.57'
$ws.Range('E35').Value = '  -6.80%  '
# Row 36
$ws.Range('D36').Value = '''0.999'
$ws.Range('E36').Value = '  -0.09%  '
# Row 37
$ws.Range('D37').Value = '''56.72'
$ws.Range('E37').Value = '  +1.20%  '
# Row 38
$ws.Range('D38').Value = '''3.34'
$ws.Range('E38').Value = '  -9.54%  '
# Row 39
$ws.Range('E39').Value = '  +0.47%  '
# Row 40
$ws.Range('D40').Value = '''33.88'
$ws.Range('E40').Value = '  +0.68%  '
# Row 41
$ws.Range('D41').Value = '''2.63'
$ws.Range('E41').Value = '  -5.57%  '
# Row 42
$ws.Range('E42').Value = '  -6.04%  '
# Row 43
$ws.Range('D43').Value = '''0.333'
$ws.Range('E43').Value = '  -2.97%  '
# Row 44
$ws.Range('D44').Value = '0.0₃0666'
$ws.Range('E44').Value = '  -5.87%  '
# Row 45
$ws.Range('E45').Value = '  -4.69%  '
# Row 46
$ws.Range('E46').Value = '  -3.93%  '
# Row 47
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '''0.128'
$ws.Range('E47').Value = '  -1.44%  '
# Row 48
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = '''2.58'
$ws.Range('E48').Value = '  -0.45%  '
# Row 49
$ws.Range('E49').Value = '  -0.03%  '
# Row 50
$ws.Range('D50').Value = '''1.35'
$ws.Range('E50').Value = '  +0.05%  '
# Row 51
$ws.Range('D51').Value = '''127.28'
$ws.Range('E51').Value = '  -3.18%  '
